$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename header suffixes: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
$headerRange = $ws.Range("A1:U1")
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2410")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2504")
    }
}

# --- 2. Turn the used range into an Excel Table (ListObject) ---
$tableRange = $ws.Range("A1:U80")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
